# Weekly update: insert the newest price record as a new row right after
# the header block's existing rows (at row 21), pushing the rest of the
# historical rows down by one and growing the used range by a row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 21; rows 21..46 shift down to 22..47.
$ws.Rows.Item(21).Insert()

# Populate the new weekly record in row 21.
$ws.Cells.Item(21, 1).Value  = 11
$ws.Cells.Item(21, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(21, 3).Value  = "Bíobío"
$ws.Cells.Item(21, 4).Value  = 44897
$ws.Cells.Item(21, 5).Value  = 8
$ws.Cells.Item(21, 6).Value  = 100112030
$ws.Cells.Item(21, 7).Value  = "Poroto granado"
$ws.Cells.Item(21, 8).Value  = "Sin especificar"
$ws.Cells.Item(21, 9).Value  = "Primera"
$ws.Cells.Item(21, 10).Value = 100
$ws.Cells.Item(21, 11).Value = 38000
$ws.Cells.Item(21, 12).Value = 40000
$ws.Cells.Item(21, 13).Value = 39000
$ws.Cells.Item(21, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(21, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(21, 16).Value = 1560
$ws.Cells.Item(21, 17).Value = 25
$ws.Cells.Item(21, 18).Value = "Hortaliza"
